$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- sheet 1 ---
$ws1.Range("F3").Value = 351
$ws1.Range("F5").Value = 7
$ws1.Range("F6").Value = 1229
$ws1.Range("F7").Value = 443
$ws1.Range("F8").Value = 100
$ws1.Range("F9").Value = 171
$ws1.Range("F10").Value = 146
$ws1.Range("F11").Value = 1038
$ws1.Range("F14").Value = 162
$ws1.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202405/HXY7pTYI1715053764601.jpeg"
$ws1.Range("F15").Value = 177
$ws1.Range("F16").Value = 1462
$ws1.Range("F17").Value = 538
$ws1.Range("F18").Value = 220
$ws1.Range("F19").Value = 335
$ws1.Range("F21").Value = 785
$ws1.Range("F22").Value = 1137
$ws1.Range("F24").Value = 1919
$ws1.Range("F25").Value = 2625
$ws1.Range("F26").Value = 1392
$ws1.Range("F27").Value = 60
$ws1.Range("F29").Value = 350
$ws1.Range("F30").Value = 400
$ws1.Range("F31").Value = 1125
$ws1.Range("F32").Value = 804
$ws1.Range("F33").Value = 1270
$ws1.Range("F34").Value = 148
$ws1.Range("F36").Value = 775
$ws1.Range("F37").Value = 552
$ws1.Range("F38").Value = 647
$ws1.Range("F39").Value = 820
$ws1.Range("F40").Value = 345
$ws1.Range("F41").Value = 235
# --- sheet 2 ---
$ws2.Range("F4").Value = 32
$ws2.Range("F8").Value = 9
$ws2.Range("F10").Value = 1
$ws2.Range("F15").Value = 600
# --- sheet 4 ---
$ws4.Range("F4").Value = 32
$ws4.Range("F6").Value = 351
$ws4.Range("F10").Value = 1229
$ws4.Range("F11").Value = 443
$ws4.Range("F12").Value = 100
$ws4.Range("F13").Value = 171
$ws4.Range("F15").Value = 146
$ws4.Range("F18").Value = 9
$ws4.Range("F19").Value = 162
$ws4.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202405/HXY7pTYI1715053764601.jpeg"
$ws4.Range("F20").Value = 177
$ws4.Range("F21").Value = 1462
$ws4.Range("F22").Value = 538
$ws4.Range("F23").Value = 220
$ws4.Range("F24").Value = 335
$ws4.Range("F26").Value = 1137
$ws4.Range("F27").Value = 2625
$ws4.Range("F29").Value = 1392
$ws4.Range("F30").Value = 60
$ws4.Range("F34").Value = 350
$ws4.Range("F35").Value = 400
$ws4.Range("F36").Value = 1125
$ws4.Range("F39").Value = 804
$ws4.Range("F40").Value = 1270
$ws4.Range("F41").Value = 775
$ws4.Range("F42").Value = 552
$ws4.Range("F43").Value = 647
$ws4.Range("F44").Value = 820
$ws4.Range("F45").Value = 346
$ws4.Range("F48").Value = 235
